$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (2 and 3) had their Fecha / Volumen / Precio values swapped.
# Capture current values first (use Value2 for reliable reads), then write the
# swapped values back.

$D2 = $ws.Range("D2").Value2
$J2 = $ws.Range("J2").Value2
$K2 = $ws.Range("K2").Value2
$L2 = $ws.Range("L2").Value2
$M2 = $ws.Range("M2").Value2
$P2 = $ws.Range("P2").Value2

$D3 = $ws.Range("D3").Value2
$J3 = $ws.Range("J3").Value2
$K3 = $ws.Range("K3").Value2
$L3 = $ws.Range("L3").Value2
$M3 = $ws.Range("M3").Value2
$P3 = $ws.Range("P3").Value2

$ws.Range("D2").Value = $D3
$ws.Range("J2").Value = $J3
$ws.Range("K2").Value = $K3
$ws.Range("L2").Value = $L3
$ws.Range("M2").Value = $M3
$ws.Range("P2").Value = $P3

$ws.Range("D3").Value = $D2
$ws.Range("J3").Value = $J2
$ws.Range("K3").Value = $K2
$ws.Range("L3").Value = $L2
$ws.Range("M3").Value = $M2
$ws.Range("P3").Value = $P2

$wb.Save()
